$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("services")

# B1 only held a blank/whitespace placeholder header with no data in the column below it.
# Clear that cell (columns C/D -- Language/Duration -- remain where they are).
$ws.Range("B1").ClearContents()

# Update the sheet's selection to match the new state.
$ws.Range("A2").Select()
